$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19 - "Projets" / "Navigation dans la liste" / "Affichage des projets suivants / précédents" (still OK)
$ws.Range("A19").Value = "Projets"
$ws.Range("B19").Value = "Navigation dans la liste"
$ws.Range("C19").Value = "Affichage des projets suivants / précédents"

# Row 20 - "Projets" / "Création d'un projet" / "Ajout du projet en haut de liste et message significatif " (KO)
$ws.Range("A20").Value = "Projets"
$ws.Range("B20").Value = "Création d'un projet"
$ws.Range("C20").Value = "Ajout du projet en haut de liste et message significatif "
$ws.Range("D20").Value = "KO"
$ws.Range("D20").Style = "Bad"

# Row 21 - "Projets" / "…" / (C stays empty) (KO)
$ws.Range("A21").Value = "Projets"
$ws.Range("B21").Value = "…"
$ws.Range("D21").Value = "KO"
$ws.Range("D21").Style = "Bad"

# Row 22 - "Équipes" / "Arrivée sur la page" / "Affichage de toutes les équipes triées par leur note" (OK)
$ws.Range("A22").Value = "Équipes"
$ws.Range("B22").Value = "Arrivée sur la page"
$ws.Range("C22").Value = "Affichage de toutes les équipes triées par leur note"

# Row 23 - "Équipes" / "Rejoindre" / "Envoie de la requête de demande d'ajout + notification" (OK)
$ws.Range("A23").Value = "Équipes"
$ws.Range("B23").Value = "Rejoindre"
$ws.Range("C23").Value = "Envoie de la requête de demande d'ajout + notification"

# Row 24 - "Équipes" / "Quitter" / "Envoie de la requête de demande d'ajout + notification" (D24 cleared entirely)
$ws.Range("A24").Value = "Équipes"
$ws.Range("B24").Value = "Quitter"
$ws.Range("C24").Value = "Envoie de la requête de demande d'ajout + notification"
$ws.Range("D24").Clear()

# Mise en page: scroll/freeze-pane anchor moved up one row and new active selection
$ws.Range("C27").Select() | Out-Null
